$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (t_period 2030): update scenario probabilities
$ws.Range("B4").Value = 0.62
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 0.08

# Row 5: relabel period to 2040, keep the new probabilities (same as 2030)
$ws.Range("A5").Value = 2040
$ws.Range("B5").Value = 0.62
$ws.Range("C5").Value = 0.3
$ws.Range("D5").Value = 0.08

# Row 6: new row for period 2050 with the old "2035" probabilities
$ws.Range("A6").Value = 2050
$ws.Range("B6").Value = 0.8099999999999999
$ws.Range("C6").Value = 0.15
$ws.Range("D6").Value = 0.04
